$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by one column (L), mirroring the existing 2020 column (K)
# so the new column gets the same style (number format/border) as K.
$ws.Range("K3").Copy() | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null
$ws.Range("L3").Value = 2020

$ws.Range("K4").Copy() | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null
$ws.Range("L4").Value = 6.18

# Leave the sheet with cell M12 selected/active, as in the authored workbook.
$ws.Range("M12").Select() | Out-Null
